$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Toggle switch  for power controlling.  "
$ws.Range("B8").Value = "Displaying  must be  on a horizontal way  as title of  temp  and speed  are on the same line  `nand the reading for each are in another line  with each  reading under its title."
$ws.Range("B10").Value = "LCD initialization  is empty with no displaying."
$ws.Range("B13").Value = "The system is preferred  to be on a PCB ."

$ws.Range("B5").Select()
